# EjemploEliminarEmp.xlsx — "Se concluye la implementacion del SSO y errores menores"
#
# The header row had its two labels swapped/renamed:
#   A1: "Empleado No."  ->  "EmpleadoNo"
#   B1: "Empresa Id"    ->  "Empresa Id"  (unchanged)
#
# The workbook is also re-saved from a newer Excel build, which naturally
# refreshes the sheet/page metadata (margins, selection, etc.) - we restore
# the recorded selection and default page orientation for fidelity.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "EmpleadoNo"
$ws.Range("B1").Value = "Empresa Id"

# Default (portrait) page orientation, as recorded by the resave.
$ws.PageSetup.Orientation = 1

# Restore the saved cursor/selection position.
$ws.Range("E15").Select()
